$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 55
$ws1.Range("F5").Value = 3992
$ws1.Range("F7").Value = 444

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 55
$ws4.Range("F5").Value = 3992
$ws4.Range("F9").Value = 444
